$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Update header row (row 1)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 values
$ws.Range("B2").Value = 16.427638310189508
$ws.Range("C2").Value = 11.112196071503888
$ws.Range("D2").Value = 13.420766835774401
$ws.Range("E2").Value = -0.25453656396425117

# Update row 3 values
$ws.Range("B3").Value = 34.135045502966477
$ws.Range("C3").Value = 3.4386750814914819
$ws.Range("D3").Value = 2.1121570451994671
$ws.Range("E3").Value = 2.879739537413883

# Update the selection range shown when the sheet is active
$ws.Range("B1:E3").Select()
